$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheets: Sheet1 -> About, BpTPEU -> BpTPEU-large, plus new BpTPEU-small
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item(1)
$wsLarge = $wb.Worksheets.Item(2)
$wsAbout.Name = "About"
$wsLarge.Name = "BpTPEU-large"

$wsSmall = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsLarge)
$wsSmall.Name = "BpTPEU-small"

# ---------------------------------------------------------------------------
# About sheet: rebuild contents from scratch
# ---------------------------------------------------------------------------
$wsAbout.Cells.Clear()

$wsAbout.Range("A1").Value = "BpTPEU BTU per Large Primary Energy Unit"
$wsAbout.Range("A1").Font.Bold = $true

$wsAbout.Range("A2").Value = "BpTPEU BTU per Small Primary Energy Unit"
$wsAbout.Range("A2").Font.Bold = $true

$wsAbout.Range("A4").Value = "Source:"
$wsAbout.Range("A4").Font.Bold = $true
$wsAbout.Range("B4").Value = "none needed"

$wsAbout.Range("B5").HorizontalAlignment = -4131  # xlLeft

$wsAbout.Range("A9").Value = "Notes"
$wsAbout.Range("A9").Font.Bold = $true

$wsAbout.Range("A10").Value = "For the U.S.:"
$wsAbout.Range("A11").Value = "The large primary energy output unit (used in totals graphs) is: quadrillion BTU"
$wsAbout.Range("A12").Value = "The small primary energy output unit (used in energy intensity per unit GDP graphs) is: thousand BTU"

$wsAbout.PageSetup.Orientation = 1  # xlPortrait

# ---------------------------------------------------------------------------
# BpTPEU-large sheet: update the unit label + its alignment
# ---------------------------------------------------------------------------
$wsLarge.Range("B1").Value = "large primary energy output unit"
$wsLarge.Range("B1").HorizontalAlignment = -4131  # xlLeft

# ---------------------------------------------------------------------------
# BpTPEU-small sheet: mirror BpTPEU-large, but with thousand BTU (10^3)
# ---------------------------------------------------------------------------
$wsSmall.Tab.ThemeColorSchemeIndex = 3

$wsSmall.Columns.Item(2).ColumnWidth = 12

$wsSmall.Range("B1").Value = "small primary energy output unit"
$wsSmall.Range("B1").HorizontalAlignment = -4131  # xlLeft

$wsSmall.Range("A2").Value = "BTU"
$wsSmall.Range("B2").Formula = "=10^3"

$wsSmall.Range("B9").Value = ""
$wsSmall.Range("B9").NumberFormat = "0.00E+00"
$wsSmall.Range("B10").Value = ""
$wsSmall.Range("B10").NumberFormat = "0.00E+00"

Write-Host "done"
